$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4
$data[0,0] = 'Bitcoin'
$data[0,1] = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$data[0,2] = '''36.805.41'
$data[0,3] = '  -0.83%  '
$data[1,0] = 'Ethereum'
$data[1,1] = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$data[1,2] = '''2.043.18'
$data[1,3] = '  -0.47%  '
$data[2,0] = 'TetherUSD'
$data[2,1] = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$data[2,2] = '''1.00'
$data[2,3] = '  +0.03%  '
$data[3,0] = 'BNB'
$data[3,1] = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$data[3,2] = '''244.59'
$data[3,3] = '  -2.02%  '
$data[4,0] = 'XRP'
$data[4,1] = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$data[4,2] = '''0.652'
$data[4,3] = '  -2.53%  '
$data[5,0] = 'USDC'
$data[5,1] = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$data[5,2] = '''1.00'
$data[5,3] = '  +0.03%  '
$data[6,0] = 'Solana'
$data[6,1] = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$data[6,2] = '''56.90'
$data[6,3] = '  -6.19%  '
$data[7,0] = 'OKB'
$data[7,1] = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$data[7,2] = '''58.33'
$data[7,3] = '  -3.53%  '
$data[8,0] = 'Cardano'
$data[8,1] = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$data[8,2] = '''0.367'
$data[8,3] = '  -5.49%  '
$data[9,0] = 'Dogecoin'
$data[9,1] = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$data[9,2] = '''0.0770'
$data[9,3] = '  -3.18%  '
$data[10,0] = 'TRON'
$data[10,1] = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$data[10,2] = '''0.110'
$data[10,3] = '  +1.32%  '
$data[11,0] = 'Chainlink'
$data[11,1] = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$data[11,2] = '''15.03'
$data[11,3] = '  -6.87%  '
$data[12,0] = 'Polygon'
$data[12,1] = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$data[12,2] = '''0.865'
$data[12,3] = '  +3.17%  '
$data[13,0] = 'WrappedliquidstakedEther2.0'
$data[13,1] = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$data[13,2] = '''2.340.95'
$data[13,3] = '  -0.44%  '
$data[14,0] = 'Polkadot'
$data[14,1] = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$data[14,2] = '''5.59'
$data[14,3] = '  -4.15%  '
$data[15,0] = 'WrappedEther'
$data[15,1] = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$data[15,2] = '''2.022.96'
$data[15,3] = '  -1.56%  '
$data[16,0] = 'Avalanche'
$data[16,1] = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$data[16,2] = '''17.83'
$data[16,3] = '  -3.60%  '
$data[17,0] = 'WrappedBTC'
$data[17,1] = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$data[17,2] = '''36.756.02'
$data[17,3] = '  -0.98%  '
$data[18,0] = 'Litecoin'
$data[18,1] = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$data[18,2] = '''73.05'
$data[18,3] = '  -3.79%  '
$data[19,0] = 'ShibaInu'
$data[19,1] = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$data[19,2] = '''0.0₃0882'
$data[19,3] = '  -2.76%  '
$data[20,0] = 'Uniswap'
$data[20,1] = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$data[20,2] = '''5.35'
$data[20,3] = '  -1.12%  '
$data[21,0] = 'BitcoinCash'
$data[21,1] = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$data[21,2] = '''235.42'
$data[21,3] = '  -1.22%  '
$data[22,0] = 'Dai'
$data[22,1] = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$data[22,2] = '''1.00'
$data[22,3] = '  +0.00%  '
$data[23,0] = 'Toncoin'
$data[23,1] = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$data[23,2] = '''2.44'
$data[23,3] = '  +0.48%  '
$data[24,0] = 'Cosmos'
$data[24,1] = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$data[24,2] = '''10.20'
$data[24,3] = '  +7.53%  '
$data[25,0] = 'PancakeSwap'
$data[25,1] = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$data[25,2] = '''2.16'
$data[25,3] = '  -2.13%  '
$data[26,0] = 'Monero'
$data[26,1] = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$data[26,2] = '''167.96'
$data[26,3] = '  -0.86%  '
$data[27,0] = 'EthereumClassic'
$data[27,1] = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$data[27,2] = '''19.80'
$data[27,3] = '  -2.11%  '
$data[28,0] = 'Filecoin'
$data[28,1] = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$data[28,2] = '''5.45'
$data[28,3] = '  +13.21%  '
$data[29,0] = 'Stellar'
$data[29,1] = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$data[29,2] = '''0.123'
$data[29,3] = '  -2.56%  '
$data[30,0] = 'ImmutableX'
$data[30,1] = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$data[30,2] = '''1.09'
$data[30,3] = '  -4.15%  '
$data[31,0] = 'InternetComputer(DFINITY)'
$data[31,1] = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$data[31,2] = '''4.67'
$data[31,3] = '  +1.83%  '
$data[32,0] = 'Hedera'
$data[32,1] = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$data[32,2] = '''0.0607'
$data[32,3] = '  -4.25%  '
$data[33,0] = 'BinanceUSD'
$data[33,1] = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$data[33,2] = '''1.00'
$data[33,3] = '  +0.13%  '
$data[34,0] = 'LidoDAOToken'
$data[34,1] = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$data[34,2] = '''2.29'
$data[34,3] = '  +2.16%  '
$data[35,0] = 'WEMIXToken'
$data[35,1] = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$data[35,2] = '''1.83'
$data[35,3] = '  +4.37%  '
$data[36,0] = 'Kaspa'
$data[36,1] = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$data[36,2] = '''0.0811'
$data[36,3] = '  -9.26%  '
$data[37,0] = 'TrustWalletToken'
$data[37,1] = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$data[37,2] = '''1.30'
$data[37,3] = '  -3.87%  '
$data[38,0] = 'THORChain'
$data[38,1] = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$data[38,2] = '''5.08'
$data[38,3] = '  -4.53%  '
$data[39,0] = 'HuobiToken'
$data[39,1] = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$data[39,2] = '''3.03'
$data[39,3] = '  -5.82%  '
$data[40,0] = 'VeChain'
$data[40,1] = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$data[40,2] = '''0.0221'
$data[40,3] = '  -1.54%  '
$data[41,0] = 'ARBITRUM'
$data[41,1] = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$data[41,2] = '''1.13'
$data[41,3] = '  -0.90%  '
$data[42,0] = 'Cronos'
$data[42,1] = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$data[42,2] = '''0.0937'
$data[42,3] = '  -14.42%  '
$data[43,0] = 'Aave'
$data[43,1] = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$data[43,2] = '''95.98'
$data[43,3] = '  -1.57%  '
$data[44,0] = 'InjectiveProtocol'
$data[44,1] = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$data[44,2] = '''16.74'
$data[44,3] = '  -5.43%  '
$data[45,0] = 'Maker'
$data[45,1] = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$data[45,2] = '''1.297.95'
$data[45,3] = '  -0.02%  '
$data[46,0] = 'RenderToken'
$data[46,1] = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$data[46,2] = '''2.34'
$data[46,3] = '  -7.34%  '
$data[47,0] = 'MXToken'
$data[47,1] = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$data[47,2] = '''2.85'
$data[47,3] = '  -1.01%  '
$data[48,0] = 'FraxShare'
$data[48,1] = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$data[48,2] = '''6.69'
$data[48,3] = '  -2.96%  '
$data[49,0] = 'RocketPoolETH'
$data[49,1] = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$data[49,2] = '''2.225.85'
$data[49,3] = '  -0.62%  '

$ws.Range("B2:E51").Value = $data

# Strip any incidental number-format styling Excel applied when the
# apostrophe-prefix forced numeric-looking Price values to remain text,
# so the cell style matches the original (unstyled) data cells.
$ws.Range("B2:E51").Style = "Normal"
